$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 880.6
$ws.Range("I32").Value = 900
$ws.Range("K32").Value = 900
$ws.Range("M32").Value = -574
$ws.Range("H88").Value = 5803.5
$ws.Range("I88").Value = 5376.625
$ws.Range("J88").Value = 6047.4287
$ws.Range("K88").Value = 5376.625
$ws.Range("L88").Value = 6047.4287
$ws.Range("M88").Value = -4970.625
$ws.Range("N88").Value = -6859.4287
$ws.Range("H91").Value = 5803.5
$ws.Range("I91").Value = 5376.625
$ws.Range("J91").Value = 6047.4287
$ws.Range("K91").Value = 5376.625
$ws.Range("L91").Value = 6047.4287
$ws.Range("M91").Value = -3972.625
$ws.Range("N91").Value = -8855.4287
$ws.Range("H129").Value = 6098594
$ws.Range("I129").Value = 35715364
$ws.Range("J129").Value = 1023.9706
$ws.Range("K129").Value = 107146092
$ws.Range("L129").Value = 3071.9118
$ws.Range("M129").Value = -107141092
$ws.Range("N129").Value = -13071.9118
$ws.Range("H137").Value = 4171880.5
$ws.Range("I137").Value = 5268428
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 15805284
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -15802734
$ws.Range("N137").Value = -20100

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4262.5713
$ws.Range("I61").Value = 962.44446
$ws.Range("J61").Value = 10202.8
$ws.Range("K61").Value = 962.44446
$ws.Range("L61").Value = 10202.8
$ws.Range("M61").Value = -750.44446
$ws.Range("N61").Value = -10626.8
$ws.Range("H122").Value = 2597.7334
$ws.Range("I122").Value = 2069
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6207
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3757
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 2846.8572
$ws.Range("I132").Value = 2689.25
$ws.Range("J132").Value = 3351.2
$ws.Range("K132").Value = 8067.75
$ws.Range("L132").Value = 10053.6
$ws.Range("M132").Value = -5537.75
$ws.Range("N132").Value = -15113.6
$ws.Range("H136").Value = 4262.5713
$ws.Range("I136").Value = 962.44446
$ws.Range("J136").Value = 10202.8
$ws.Range("K136").Value = 2887.33338
$ws.Range("L136").Value = 30608.4
$ws.Range("M136").Value = -337.33338
$ws.Range("N136").Value = -35708.39999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2771.484
$ws.Range("I134").Value = 1918.12
$ws.Range("J134").Value = 6327.1665
$ws.Range("K134").Value = 5754.36
$ws.Range("L134").Value = 18981.4995
$ws.Range("M134").Value = -3219.36
$ws.Range("N134").Value = -24051.4995

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23814306
$ws.Range("J58").Value = 45460172
$ws.Range("L58").Value = 45460172
$ws.Range("N58").Value = -45460578
$ws.Range("H136").Value = 23814306
$ws.Range("J136").Value = 45460172
$ws.Range("L136").Value = 136380516
$ws.Range("N136").Value = -136385616

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 13194.2
$ws.Range("J63").Value = 15741.25
$ws.Range("L63").Value = 47223.75
$ws.Range("N63").Value = -48721.75
$ws.Range("H66").Value = 13194.2
$ws.Range("J66").Value = 15741.25
$ws.Range("L66").Value = 141671.25
$ws.Range("N66").Value = -149159.25
$ws.Range("H113").Value = 3031159.2
$ws.Range("I113").Value = 8333919
$ws.Range("J113").Value = 1011
$ws.Range("K113").Value = 25001757
$ws.Range("L113").Value = 3033
$ws.Range("M113").Value = -24999587
$ws.Range("N113").Value = -7373

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4505.684
$ws.Range("I122").Value = 3316
$ws.Range("K122").Value = 9948
$ws.Range("M122").Value = -7498
$ws.Range("H126").Value = 3504.9473
$ws.Range("I126").Value = 1193.3334
$ws.Range("J126").Value = 3938.375
$ws.Range("K126").Value = 3580.0002
$ws.Range("L126").Value = 11815.125
$ws.Range("M126").Value = -1110.0002
$ws.Range("N126").Value = -16755.125
$ws.Range("H132").Value = 4490.759
$ws.Range("I132").Value = 5014.625
$ws.Range("J132").Value = 3846
$ws.Range("K132").Value = 15043.875
$ws.Range("L132").Value = 11538
$ws.Range("M132").Value = -12513.875
$ws.Range("N132").Value = -16598

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1850.421
$ws.Range("I100").Value = 1210.7273
$ws.Range("K100").Value = 1210.7273
$ws.Range("M100").Value = -669.7273
$ws.Range("H122").Value = 2564.1738
$ws.Range("I122").Value = 2239.5293
$ws.Range("J122").Value = 3484
$ws.Range("K122").Value = 6718.5879
$ws.Range("L122").Value = 10452
$ws.Range("M122").Value = -4268.5879
$ws.Range("N122").Value = -15352
$ws.Range("H132").Value = 2316.0857
$ws.Range("I132").Value = 1700.3334
$ws.Range("J132").Value = 2777.9
$ws.Range("K132").Value = 5101.0002
$ws.Range("L132").Value = 8333.700000000001
$ws.Range("M132").Value = -2571.0002
$ws.Range("N132").Value = -13393.7
$ws.Range("H136").Value = 6253437
$ws.Range("I136").Value = 20003618
$ws.Range("J136").Value = 3354.5454
$ws.Range("K136").Value = 60010854
$ws.Range("L136").Value = 10063.6362
$ws.Range("M136").Value = -60008304
$ws.Range("N136").Value = -15163.6362

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 30328.572
$ws.Range("J82").Value = 30328.572
$ws.Range("L82").Value = 30328.572
$ws.Range("N82").Value = -31094.572
$ws.Range("H85").Value = 30328.572
$ws.Range("J85").Value = 30328.572
$ws.Range("L85").Value = 30328.572
$ws.Range("N85").Value = -32980.572
$ws.Range("H92").Value = 39250
$ws.Range("J92").Value = 39250
$ws.Range("L92").Value = 39250
$ws.Range("N92").Value = -44242
$ws.Range("H96").Value = 1332.5385
$ws.Range("I96").Value = 1238.9
$ws.Range("J96").Value = 1644.6666
$ws.Range("K96").Value = 1238.9
$ws.Range("L96").Value = 1644.6666
$ws.Range("M96").Value = 134.0999999999999
$ws.Range("N96").Value = -4390.6666
$ws.Range("H100").Value = 496.66666
$ws.Range("I100").Value = 445
$ws.Range("K100").Value = 890
$ws.Range("M100").Value = -349
$ws.Range("H122").Value = 478352.56
$ws.Range("I122").Value = 770984.25
$ws.Range("K122").Value = 2312952.75
$ws.Range("M122").Value = -2310502.75
$ws.Range("H135").Value = 83571.336
$ws.Range("J135").Value = 83571.336
$ws.Range("L135").Value = 83571.336
$ws.Range("N135").Value = -93711.336
